# Scheduled-runner update: refresh Marketboard/Leve profit figures (H,I,J,K,L,M,N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 3534.4138
$ws.Range("I40").Value = 1199.5
$ws.Range("J40").Value = 3908
$ws.Range("K40").Value = 1199.5
$ws.Range("L40").Value = 3908
$ws.Range("M40").Value = -1024.5
$ws.Range("N40").Value = -4258
# Row 51
$ws.Range("H51").Value = 3045.516
$ws.Range("I51").Value = 3574.3333
$ws.Range("K51").Value = 3574.3333
$ws.Range("M51").Value = -3090.3333
# Row 62
$ws.Range("H62").Value = 7523.067
$ws.Range("I62").Value = 6638.0835
$ws.Range("K62").Value = 6638.0835
$ws.Range("M62").Value = -6014.0835
# Row 65
$ws.Range("H65").Value = 7523.067
$ws.Range("I65").Value = 6638.0835
$ws.Range("K65").Value = 33190.4175
$ws.Range("M65").Value = -30070.4175
# Row 74
$ws.Range("H74").Value = 11655.444
$ws.Range("I74").Value = 15725
$ws.Range("J74").Value = 8399.799999999999
$ws.Range("K74").Value = 15725
$ws.Range("L74").Value = 8399.799999999999
$ws.Range("M74").Value = -14789
$ws.Range("N74").Value = -10271.8
# Row 76
$ws.Range("H76").Value = 5371
$ws.Range("I76").Value = 4875.75
$ws.Range("K76").Value = 4875.75
$ws.Range("M76").Value = -4560.75
# Row 77
$ws.Range("H77").Value = 11655.444
$ws.Range("I77").Value = 15725
$ws.Range("J77").Value = 8399.799999999999
$ws.Range("K77").Value = 78625
$ws.Range("L77").Value = 41999
$ws.Range("M77").Value = -73945
$ws.Range("N77").Value = -51359
# Row 79
$ws.Range("H79").Value = 5371
$ws.Range("I79").Value = 4875.75
$ws.Range("K79").Value = 4875.75
$ws.Range("M79").Value = -3783.75
# Row 93
$ws.Range("H93").Value = 84996.336
$ws.Range("J93").Value = 84996.336
$ws.Range("L93").Value = 84996.336
$ws.Range("N93").Value = -89988.336
# Row 131
$ws.Range("H131").Value = 3119.6
$ws.Range("I131").Value = 2649.5
$ws.Range("J131").Value = 3433
$ws.Range("K131").Value = 7948.5
$ws.Range("L131").Value = 10299
$ws.Range("M131").Value = -2908.5
$ws.Range("N131").Value = -20379
# Row 141
$ws.Range("H141").Value = 4701.1665
$ws.Range("I141").Value = 4333.4
$ws.Range("K141").Value = 13000.2
$ws.Range("M141").Value = -7820.199999999999
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 25001242
$ws.Range("I32").Value = 29412726
$ws.Range("J32").Value = 2834.6667
$ws.Range("K32").Value = 29412726
$ws.Range("L32").Value = 2834.6667
$ws.Range("M32").Value = -29412439
$ws.Range("N32").Value = -3408.6667
# Row 110
$ws.Range("H110").Value = 2135
$ws.Range("I110").Value = 987.7778
$ws.Range("K110").Value = 987.7778
$ws.Range("M110").Value = 1057.2222
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 494.5357
$ws.Range("I94").Value = 573.5
$ws.Range("J94").Value = 205
$ws.Range("K94").Value = 573.5
$ws.Range("L94").Value = 205
$ws.Range("M94").Value = -122.5
$ws.Range("N94").Value = -1107
# Row 97
$ws.Range("H97").Value = 27374.666
$ws.Range("I97").Value = 6491.857
$ws.Range("K97").Value = 6491.857
$ws.Range("M97").Value = -5500.857
# Row 105
$ws.Range("H105").Value = 2976.4285
$ws.Range("I105").Value = 3032.8
$ws.Range("K105").Value = 3032.8
$ws.Range("M105").Value = -1285.8
# Row 134
$ws.Range("H134").Value = 3577.6538
$ws.Range("I134").Value = 2343.889
$ws.Range("J134").Value = 6353.625
$ws.Range("K134").Value = 7031.667
$ws.Range("L134").Value = 19060.875
$ws.Range("M134").Value = -4496.667
$ws.Range("N134").Value = -24130.875
$ws = $wb.Worksheets.Item("CRP")
# Row 59
$ws.Range("H59").Value = 80000
$ws.Range("J59").Value = 80000
$ws.Range("L59").Value = 80000
$ws.Range("N59").Value = -82290
# Row 132
$ws.Range("H132").Value = 4901.9614
$ws.Range("I132").Value = 4901.9614
$ws.Range("K132").Value = 14705.8842
$ws.Range("M132").Value = -12175.8842
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 558.4
$ws.Range("I5").Value = 520.4
$ws.Range("J5").Value = 596.4
$ws.Range("K5").Value = 1561.2
$ws.Range("L5").Value = 1789.2
$ws.Range("M5").Value = -1449.2
$ws.Range("N5").Value = -2013.2
# Row 68
$ws.Range("H68").Value = 1721.1515
$ws.Range("I68").Value = 1328.5
$ws.Range("J68").Value = 2010.4736
$ws.Range("K68").Value = 3985.5
$ws.Range("L68").Value = 6031.4208
$ws.Range("M68").Value = -3174.5
$ws.Range("N68").Value = -7653.4208
# Row 71
$ws.Range("H71").Value = 1721.1515
$ws.Range("I71").Value = 1328.5
$ws.Range("J71").Value = 2010.4736
$ws.Range("K71").Value = 11956.5
$ws.Range("L71").Value = 18094.2624
$ws.Range("M71").Value = -7900.5
$ws.Range("N71").Value = -26206.2624
# Row 131
$ws.Range("H131").Value = 1930.9642
$ws.Range("J131").Value = 1745.3877
$ws.Range("L131").Value = 5236.1631
$ws.Range("N131").Value = -15316.1631
# Row 135
$ws.Range("H135").Value = 558.4
$ws.Range("I135").Value = 520.4
$ws.Range("J135").Value = 596.4
$ws.Range("K135").Value = 4683.599999999999
$ws.Range("L135").Value = 5367.599999999999
$ws.Range("M135").Value = -2148.599999999999
$ws.Range("N135").Value = -10437.6
# Row 140
$ws.Range("H140").Value = 920.5625
$ws.Range("J140").Value = 1200
$ws.Range("L140").Value = 3600
$ws.Range("N140").Value = -13960
# Row 141
$ws.Range("H141").Value = 12250
$ws.Range("I141").Value = 12250
$ws.Range("K141").Value = 36750
$ws.Range("M141").Value = -31570
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3125
$ws.Range("I80").Value = 1500
$ws.Range("K80").Value = 1500
$ws.Range("M80").Value = -502
# Row 83
$ws.Range("H83").Value = 3125
$ws.Range("I83").Value = 1500
$ws.Range("K83").Value = 7500
$ws.Range("M83").Value = -2508
# Row 102
$ws.Range("H102").Value = 2114.8147
$ws.Range("I102").Value = 2085.4583
$ws.Range("J102").Value = 2349.6667
$ws.Range("K102").Value = 2085.4583
$ws.Range("L102").Value = 2349.6667
$ws.Range("M102").Value = -463.4582999999998
$ws.Range("N102").Value = -5593.6667
# Row 104
$ws.Range("H104").Value = 144199.75
$ws.Range("J104").Value = 144199.75
$ws.Range("L104").Value = 144199.75
$ws.Range("N104").Value = -151187.75
# Row 126
$ws.Range("H126").Value = 4068.7693
$ws.Range("I126").Value = 3999.75
$ws.Range("J126").Value = 4099.4443
$ws.Range("K126").Value = 11999.25
$ws.Range("L126").Value = 12298.3329
$ws.Range("M126").Value = -9529.25
$ws.Range("N126").Value = -17238.3329
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 7660.609
$ws.Range("J46").Value = 8214.048000000001
$ws.Range("L46").Value = 8214.048000000001
$ws.Range("N46").Value = -8590.048000000001
# Row 82
$ws.Range("H82").Value = 4187.5
$ws.Range("I82").Value = 5000
$ws.Range("J82").Value = 3916.6667
$ws.Range("K82").Value = 5000
$ws.Range("L82").Value = 3916.6667
$ws.Range("M82").Value = -4639
$ws.Range("N82").Value = -4638.6667
# Row 85
$ws.Range("H85").Value = 4187.5
$ws.Range("I85").Value = 5000
$ws.Range("J85").Value = 3916.6667
$ws.Range("K85").Value = 5000
$ws.Range("L85").Value = 3916.6667
$ws.Range("M85").Value = -3752
$ws.Range("N85").Value = -6412.6667
# Row 93
$ws.Range("H93").Value = 71430510
$ws.Range("I93").Value = 125001736
$ws.Range("J93").Value = 2216.6667
$ws.Range("K93").Value = 125001736
$ws.Range("L93").Value = 2216.6667
$ws.Range("M93").Value = -125000488
$ws.Range("N93").Value = -4712.6667
# Row 136
$ws.Range("H136").Value = 20051.727
$ws.Range("J136").Value = 4553.9
$ws.Range("L136").Value = 13661.7
$ws.Range("N136").Value = -18761.7
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 5301.773
$ws.Range("I81").Value = 1386.5834
$ws.Range("K81").Value = 2773.1668
$ws.Range("M81").Value = -1712.1668
# Row 84
$ws.Range("H84").Value = 5301.773
$ws.Range("I84").Value = 1386.5834
$ws.Range("K84").Value = 13865.834
$ws.Range("M84").Value = -8561.833999999999
# Row 100
$ws.Range("H100").Value = 575.8
$ws.Range("I100").Value = 561.5833
$ws.Range("K100").Value = 1123.1666
$ws.Range("M100").Value = -582.1666
# Row 137
$ws.Range("H137").Value = 118995
$ws.Range("J137").Value = 118995
$ws.Range("L137").Value = 118995
$ws.Range("N137").Value = -129195
